# Append the new daily profit-data row (row 82) to the bottom of the sheet,
# mirroring the existing rows' layout: col A = date text, cols B:L = numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 82

# Column A holds dates stored as plain text (e.g. "02/14/2026"), just like
# every other row in the sheet. Pre-format the cell as Text so Excel does not
# silently convert the date-looking string into a real date serial number.
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "02/14/2026"

$ws.Range("B$newRow").Value = 9731
$ws.Range("C$newRow").Value = 0.2354351560586131
$ws.Range("D$newRow").Value = 0.7645648439413869
$ws.Range("E$newRow").Value = -294.16
$ws.Range("F$newRow").Value = -33.24
$ws.Range("G$newRow").Value = -23501.78
$ws.Range("H$newRow").Value = -75.95
$ws.Range("I$newRow").Value = -1061.83
$ws.Range("J$newRow").Value = -31.67
$ws.Range("K$newRow").Value = -24563.61
$ws.Range("L$newRow").Value = -71.63
